# Update figures and figure data
#
# fig3.xlsx gains a second worksheet ("Sheet1") that holds the data behind a
# new "flfp" (female labour-force participation) breakdown, right after the
# existing "fig3" sheet. The new sheet becomes the active tab/sheet; the old
# "fig3" sheet loses its tab-selected flag and both sheets end up with a
# C1:F1 header selection.

$wb   = $excel.ActiveWorkbook
$fig3 = $wb.Worksheets.Item("fig3")

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet immediately after "fig3" and rename it.
# ---------------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $fig3)
$new.Name = "Sheet1"

# ---------------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------------
$new.Range("A1").Value = "flfp"
$new.Range("B1").Value = "marst"
$new.Range("C1").Value = "One `$ Manager"
$new.Range("C1").NumberFormat = "0.00"
$new.Range("D1").Value = "Manage `$ Together"
$new.Range("E1").Value = "Keep Some `$ Separate"
$new.Range("F1").Value = "Keep All `$ Separate"
$new.Range("F1").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 3. Data rows - age index (A), marital status (B) and four share columns
#    (C:F). Rows 2-14 are "Cohab", rows 15-27 are "Married".
# ---------------------------------------------------------------------------
$data = @(
  @(39, "Cohab",   0.4878579,  0.35235632, 0.09203948, 0.0677463),
  @(40, "Cohab",   0.43218853, 0.36890349, 0.11479871, 0.08410927),
  @(41, "Cohab",   0.37675539, 0.37992325, 0.14072254, 0.10259882),
  @(42, "Cohab",   0.32290244, 0.38468476, 0.1694884,  0.12292441),
  @(43, "Cohab",   0.27191763, 0.38285039, 0.20057774, 0.14465424),
  @(44, "Cohab",   0.22491229, 0.37452496, 0.23331306, 0.16724968),
  @(45, "Cohab",   0.1827223,  0.36024382, 0.26691791, 0.19011598),
  @(46, "Cohab",   0.1458493,  0.34090199, 0.30058911, 0.2126596),
  @(47, "Cohab",   0.11445014, 0.31764006, 0.33356837, 0.23434143),
  @(48, "Cohab",   0.08837115, 0.29171057, 0.36520113, 0.25471716),
  @(49, "Cohab",   0.06721506, 0.26434967, 0.39497513, 0.27346014),
  @(50, "Cohab",   0.05042361, 0.23667357, 0.42253622, 0.29036661),
  @(51, "Cohab",   0.0373596,  0.20961039, 0.44768382, 0.3053462),
  @(39, "Married", 0.46413617, 0.41521345, 0.06820083, 0.05244955),
  @(40, "Married", 0.41642708, 0.44798033, 0.07863813, 0.05695446),
  @(41, "Married", 0.37004248, 0.47884846, 0.08983886, 0.06127021),
  @(42, "Married", 0.32570585, 0.50724239, 0.10172946, 0.06532229),
  @(43, "Married", 0.28402416, 0.53269674, 0.11423004, 0.06904907),
  @(44, "Married", 0.24546096, 0.55487485, 0.12726002, 0.07240417),
  @(45, "Married", 0.21032428, 0.57357517, 0.14074327, 0.07535728),
  @(46, "Married", 0.17876857, 0.58872583, 0.15461201, 0.07789358),
  @(47, "Married", 0.15080865, 0.60036992, 0.16880939, 0.08001204),
  @(48, "Married", 0.12634193, 0.60864449, 0.18329059, 0.08172299),
  @(49, "Married", 0.10517524, 0.61375674, 0.19802264, 0.08304538),
  @(50, "Married", 0.08705231, 0.61596037, 0.21298329, 0.08400403),
  @(51, "Married", 0.07167949, 0.61553418, 0.22815912, 0.08462721)
)

# Write the C:F (percentage) columns first so the two new cell styles are
# appended to styles.xml's cellXfs in the same order as the target file
# (numFmt "0.00" + center alignment, then default numFmt + center alignment).
$r = 2
foreach ($row in $data) {
    $new.Range("C$r").Value = $row[2]
    $new.Range("D$r").Value = $row[3]
    $new.Range("E$r").Value = $row[4]
    $new.Range("F$r").Value = $row[5]
    $r = $r + 1
}
$new.Range("C2:F27").NumberFormat = "0.00"
$new.Range("C2:F27").HorizontalAlignment = -4108

# Now write column A/B and center column A.
$r = 2
foreach ($row in $data) {
    $new.Range("A$r").Value = $row[0]
    $new.Range("B$r").Value = $row[1]
    $r = $r + 1
}
$new.Range("A2:A27").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Recreate the leftover sortState (the sheet was last sorted by the
#    "marst" column, B2:B27 over A2:F27).
# ---------------------------------------------------------------------------
$sortRange = $new.Range("A2:F27")
$sortKey   = $new.Range("B2:B27")
$new.Sort.SortFields.Clear()
$new.Sort.SortFields.Add($sortKey)
$new.Sort.SetRange($sortRange)
$new.Sort.Apply()

# ---------------------------------------------------------------------------
# 5. Selection / active-sheet state.
# ---------------------------------------------------------------------------
$fig3.Range("C1:F1").Select()

$new.Range("C1:F1").Select()
$new.Activate()
